$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.817.48'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -6.88%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.537.51'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.01%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '298.36'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.62%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '93.98'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.73%  '

# Row 7
$ws.Range('E7').Value = '  -3.95%  '

# Row 8
$ws.Range('E8').Value = '  +0.07%  '

# Row 9
$ws.Range('E9').Value = '  -5.28%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.00'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -7.37%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0802'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.50%  '

# Row 12
$ws.Range('E12').Value = '  -4.59%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.113'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.88%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.932.07'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.63%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.521.47'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.68%  '

# Row 16
$ws.Range('E16').Value = '  -5.06%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.05'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -5.12%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.896.08'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -6.69%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.87'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.58%  '

# Row 20
$ws.Range('E20').Value = '  -3.23%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.57'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.03%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.58'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.89%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '255.72'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -9.98%  '

# Row 24
$ws.Range('E24').Value = '  -4.27%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -7.23%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '28.96'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.84%  '

# Row 27
$ws.Range('E27').Value = '  +0.06%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.98'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.42%  '

# Row 29
$ws.Range('E29').Value = '  -4.43%  '

# Row 30
$ws.Range('E30').Value = '  -3.92%  '

# Row 31
$ws.Range('E31').Value = '  -5.52%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '152.53'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.88%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.75'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.88%  '

# Row 34
$ws.Range('E34').Value = '  -6.18%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.38'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -6.46%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0794'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.72%  '

# Row 37
$ws.Range('E37').Value = '  -5.72%  '

# Row 38
$ws.Range('E38').Value = '  -2.70%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '16.84'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.71%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '23.28'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.89%  '

# Row 41
$ws.Range('E41').Value = '  -3.52%  '

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0310'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.26%  '

# Row 43
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.88'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.67%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.084.49'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.05%  '

# Row 45
$ws.Range('E45').Value = '  +0.17%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '83.65'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -11.09%  '

# Row 47
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.90'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.22%  '

# Row 48
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.59'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.43%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.787.94'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.69%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '104.24'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.68%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.65'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.65%  '
